# Mark additional checklist rows as completed ("x") in the "Buffer Days Plan"
# sheet's completion column (E), matching the rows whose tasks have now been
# documented/finished, then leave the selection where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buffer Days Plan")
$ws.Activate()

$completedRows = @(36,37,41,42,43,47,49,52,56,62,63,64,65,66,68,69,70,71,72,77,78)

foreach ($r in $completedRows) {
    $ws.Cells.Item($r, 5).Value = "x"
}

$ws.Range("E111").Select()
